# Fix a typo in the "Challenges" bullet list on the Lessons-Learned /
# Challenges slide (sldId 269, the 11th slide): remove the stray "it "
# from "...as it we had to use a spreadsheet..." so it reads
# "...as we had to use a spreadsheet...".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(9)   # "TextBox 14"

# Editing text re-triggers the shape's auto-fit height calculation, so
# remember the current height and restore it once we're done.
$origHeight = $shp.Height

$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)   # "Loading Box Office data on UI ..." bullet

# Rewrite the whole run's characters (rather than just the changed
# substring) so the paragraph keeps a single run instead of being split
# into multiple runs with identical formatting.
$run = $para.Characters(1, $para.Length)
$run.Text = "Loading Box Office data on UI was one of the important features of the tool in which we faced challenges as we had to use a spreadsheet that was generated via a third-party application."

# Restore the shape height. The Height property round-trips through a
# 32-bit float and truncates toward the EMU below, so nudge it up by a
# hair to land back on the exact original EMU value.
$shp.Height = $origHeight + 0.00003
